# DaySale_2025-08-23_00-00.xlsx update:
#   - a new item "CONCOR 5MG 30 TAB" was recorded, so it now shows up in the
#     "low stock" list as item #3 and every item after it shifts down by one
#     slot (the last item, سرنجات 3 سم, lands on a freshly inserted row).
#   - the generated-on timestamp moved from 9:36 AM to 9:41 AM.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. make room for the item that gets pushed off the bottom of the table ---
# Insert one row right before the current "total" row (row 14) and give it the
# same look as the row above it (row 13, the last item row).
$ws.Rows(14).Insert()
$ws.Range("A13:Q13").Copy()
$ws.Range("A14:Q14").PasteSpecial(-4122)
$ws.Rows(14).RowHeight = 25.5

$ws.Range("A14:B14").Merge()
$ws.Range("C14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("N14:O14").Merge()

# --- 2. shift the items: row 9 becomes the new item, rows 10-13 pick up the
#        values that used to sit one row above them, and the new row 14 gets
#        what used to be the last item (سرنجات 3 سم) ---

# row 9  (was item 3 - PANTOLOC)      -> new item 3 - CONCOR 5MG 30 TAB
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "CONCOR 5MG 30 TAB"
$ws.Range("H9").Value = "2:0"
$ws.Range("L9").Value = "1"
$ws.Range("N9").Value = "72.00"
$ws.Range("P9").Value = "47.5200"
$ws.Range("Q9").Value = "0:2"

# row 10 (was item 4 - PROPAMETHONE)  -> item 4 - PANTOLOC 40MG 14 TAB
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "PANTOLOC 40MG 14 TAB"
$ws.Range("H10").Value = "1:1"
$ws.Range("L10").Value = "1"
$ws.Range("N10").Value = "102.00"
$ws.Range("P10").Value = "51.0000"
$ws.Range("Q10").Value = "0:1"

# row 11 (was item 5 - QUIBRON)       -> item 5 - PROPAMETHONE TOP. CREAM. 20 GM
$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "PROPAMETHONE TOP. CREAM. 20 GM"
$ws.Range("H11").Value = "1:0"
$ws.Range("L11").Value = "1"
$ws.Range("N11").Value = "30.00"
$ws.Range("P11").Value = "30.0000"
$ws.Range("Q11").Value = "1:0"

# row 12 (was item 6 - TAVACIN)       -> item 6 - QUIBRON T/SR 300MG 100 TAB
$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "QUIBRON T/SR 300MG 100 TAB"
$ws.Range("H12").Value = "1:41"
$ws.Range("L12").Value = "1"
$ws.Range("N12").Value = "132.00"
$ws.Range("P12").Value = "6.6000"
$ws.Range("Q12").Value = "0:5"

# row 13 (was item 7 - سرنجات 3 سم)   -> item 7 - TAVACIN 500MG 5 F.C. TAB
$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "TAVACIN 500MG 5 F.C. TAB"
$ws.Range("H13").Value = "1:0"
$ws.Range("L13").Value = "1"
$ws.Range("N13").Value = "90.00"
$ws.Range("P13").Value = "90.0000"
$ws.Range("Q13").Value = "1:0"

# row 14 (new row)                    -> item 8 - سرنجات 3 سم
$ws.Range("A14").Value = 8
$ws.Range("C14").Value = "سرنجات 3 سم"
$ws.Range("H14").Value = "0:0"
$ws.Range("L14").Value = "0"
$ws.Range("N14").Value = "2.00"
$ws.Range("P14").Value = "2.0000"
$ws.Range("Q14").Value = "1:0"

# --- 3. update the total shown under the price column ---
$ws.Range("P15").Value = 278.91000000000003

# --- 4. refresh the generated-on timestamp in the footer ---
$ws.Range("A16").Value = "Saturday, 23 August, 2025 9:41 AM"
